# Auto-generated edit script: updates crypto price/volume values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.630.23"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "1.593.59"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  -1.42%  "
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.514"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.247"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0617"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  -1.76%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.66"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  -2.47%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0836"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "1.817.48"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "1.588.00"
$ws.Range("E13").Value = "  -2.02%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.02"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("E15").Value = "  -3.09%  "
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "65.16"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "26.614.21"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -1.90%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "210.16"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("E20").Value = "  +0.10%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.70"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("E22").Value = "  -2.66%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.31"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -2.80%  "
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.88"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -2.22%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "146.39"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E32").Value = "  -3.48%  "
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.671"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -11.56%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.90"
$cell.Style = $origStyle
$ws.Range("D35").Value = "1.292.28"
$ws.Range("E35").Value = "  -4.39%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -5.34%  "
$ws.Range("E38").Value = "  -3.21%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.834"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  +0.51%  "
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "63.59"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").Value = "1.729.96"
$ws.Range("E45").Value = "  -1.71%  "
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "89.54"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  -0.63%  "
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.61"
$cell.Style = $origStyle
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.832"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -5.71%  "
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0986"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("E50").Value = "  -2.36%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.51"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -2.30%  "
